$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.233.87"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "3.369.57"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.35"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.14"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D8").Value = "3.369.68"
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.472"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.48"
$ws.Range("E10").Value = "  -1.62%  "
$ws.Range("E11").Value = "  +1.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.388"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "3.947.12"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000175"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.88"
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "3.375.58"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "61.347.81"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.06"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.85"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.33"
$ws.Range("E21").Value = "  -0.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.41"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.554"
$ws.Range("E23").Value = "  -3.38%  "
$ws.Range("D24").Value = "3.511.92"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("E26").Value = "  +7.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.43"
$ws.Range("E27").Value = "  +0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.71"
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.42"
$ws.Range("E30").Value = "  -4.37%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.23"
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.159"
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.16"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.46"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.25"
$ws.Range("E36").Value = "  -5.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.81"
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.53"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.65"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0773"
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.773"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.72"
$ws.Range("E43").Value = "  +5.61%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.40"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.41"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.20"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.57"
$ws.Range("E47").Value = "  +5.87%  "
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.59"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("D50").Value = "2.344.63"
$ws.Range("E50").Value = "  +3.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0260"
$ws.Range("E51").Value = "  +0.06%  "
